$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.710.92'
$ws.Range("E2").Value = '  +3.15%  '
$ws.Range("D3").Value = '2.194.49'
$ws.Range("E3").Value = '  +0.73%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'259.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.61%  '
$ws.Range("D6").Value = "'81.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.71%  '
$ws.Range("E7").Value = '  +2.97%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'0.594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").Value = "'43.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.17%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  +3.23%  '
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").Value = '2.517.34'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = "'14.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").Value = '2.183.71'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = "'0.777"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").Value = '43.608.06'
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Value = "'70.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").Value = "'5.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").Value = "'2.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.20%  '
$ws.Range("D23").Value = "'230.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").Value = "'8.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.96%  '
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").Value = "'42.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +14.93%  '
$ws.Range("D27").Value = "'10.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("E30").Value = '  +2.54%  '
$ws.Range("D31").Value = "'173.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("D33").Value = "'0.0870"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.48%  '
$ws.Range("D34").Value = "'5.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.89%  '
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("D37").Value = "'4.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.42%  '
$ws.Range("D38").Value = "'0.0351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.12%  '
$ws.Range("D39").Value = "'13.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.96%  '
$ws.Range("E40").Value = '  +15.97%  '
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("D42").Value = "'62.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.33%  '
$ws.Range("D43").Value = "'5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.58%  '
$ws.Range("E44").Value = '  +1.46%  '
$ws.Range("D45").Value = "'100.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.54%  '
$ws.Range("D46").Value = "'0.0981"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").Value = "'8.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  +4.22%  '
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D51").Value = "'0.438"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.16%  '
